$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "-"
$ws.Range("F2").Value = "[-, Emerson-Comandos Eletricos-1B, -, Weslei-Metrologia 1-1B]"

$ws.Range("B3").Value = "-"
$ws.Range("C3").Value = "-"
$ws.Range("E3").Value = "[Emerson-Comandos Eletricos-1B, -, Weslei-Metrologia 1-1B, -]"
$ws.Range("F3").Value = "[Ernane-Desenho tecnico mecanico-1B, Ernane-Desenho tecnico mecanico-1B]"

$ws.Range("B4").Value = "-"
$ws.Range("C4").Value = "-"
$ws.Range("E4").Value = "[Emerson-Comandos Eletricos-1B, -, Weslei-Metrologia 1-1B, -]"
$ws.Range("F4").Value = "[-, Ernane-Desenho tecnico mecanico-1B]"

$ws.Range("B6").Value = "-"
$ws.Range("E6").Value = "[Aline S. M.-T. M. Metalicos-1B, Aline S. M.-T. M. Metalicos-1B]"
$ws.Range("F6").Value = "Anselmo-Gestao Intregrada"

$ws.Range("B7").Value = "-"
$ws.Range("F7").Value = "Anselmo-Gestao Intregrada"

$ws.Range("B8").Value = "-"
$ws.Range("F8").Value = "[Weslei-Metrologia 1-1B, -, Emerson-Comandos Eletricos-1B, -]"

$ws.Range("C20").Value = "-"
